$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3818831.8
$ws.Range("I40").Value = 6946989
$ws.Range("J40").Value = 1003490.2
$ws.Range("K40").Value = 6946989
$ws.Range("L40").Value = 1003490.2
$ws.Range("M40").Value = -6946814
$ws.Range("N40").Value = -1003840.2
$ws.Range("H41").Value = 149.22223
$ws.Range("I41").Value = 100.63636
$ws.Range("K41").Value = 100.63636
$ws.Range("M41").Value = 339.36364
$ws.Range("H52").Value = 83336830
$ws.Range("J52").Value = 100003590
$ws.Range("L52").Value = 300010770
$ws.Range("N52").Value = -300011090
$ws.Range("H53").Value = 18519088
$ws.Range("I53").Value = 35714436
$ws.Range("K53").Value = 35714436
$ws.Range("M53").Value = -35713799
$ws.Range("H62").Value = 15634400
$ws.Range("I62").Value = 22738672
$ws.Range("J62").Value = 5001.2
$ws.Range("K62").Value = 22738672
$ws.Range("L62").Value = 5001.2
$ws.Range("M62").Value = -22738048
$ws.Range("N62").Value = -6249.2
$ws.Range("H65").Value = 15634400
$ws.Range("I65").Value = 22738672
$ws.Range("J65").Value = 5001.2
$ws.Range("K65").Value = 113693360
$ws.Range("L65").Value = 25006
$ws.Range("M65").Value = -113690240
$ws.Range("N65").Value = -31246
$ws.Range("H114").Value = 31850
$ws.Range("J114").Value = 31850
$ws.Range("L114").Value = 31850
$ws.Range("N114").Value = -40528
$ws.Range("H129").Value = 1188.6538
$ws.Range("J129").Value = 1016.8333
$ws.Range("L129").Value = 3050.4999
$ws.Range("N129").Value = -13050.4999
$ws.Range("H137").Value = 1289.5918
$ws.Range("I137").Value = 1048.4634
$ws.Range("J137").Value = 2525.375
$ws.Range("K137").Value = 3145.3902
$ws.Range("L137").Value = 7576.125
$ws.Range("M137").Value = -595.3902000000003
$ws.Range("N137").Value = -12676.125
$ws.Range("H138").Value = 1580.0596
$ws.Range("I138").Value = 664.55554
$ws.Range("J138").Value = 3227.9666
$ws.Range("K138").Value = 1993.66662
$ws.Range("L138").Value = 9683.899800000001
$ws.Range("M138").Value = 3146.33338
$ws.Range("N138").Value = -19963.8998

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 5628.5713
$ws.Range("I21").Value = 1400
$ws.Range("J21").Value = 6333.3335
$ws.Range("K21").Value = 1400
$ws.Range("L21").Value = 6333.3335
$ws.Range("M21").Value = -1026
$ws.Range("N21").Value = -7081.3335
$ws.Range("H32").Value = 12349733
$ws.Range("I32").Value = 3732.3896
$ws.Range("K32").Value = 3732.3896
$ws.Range("M32").Value = -3445.3896
$ws.Range("H45").Value = 42823.75
$ws.Range("I45").Value = 63433.938
$ws.Range("J45").Value = 1603.375
$ws.Range("K45").Value = 63433.938
$ws.Range("L45").Value = 1603.375
$ws.Range("M45").Value = -63056.938
$ws.Range("N45").Value = -2357.375
$ws.Range("H61").Value = 1077.973
$ws.Range("I61").Value = 1081.3636
$ws.Range("J61").Value = 1050
$ws.Range("K61").Value = 1081.3636
$ws.Range("L61").Value = 1050
$ws.Range("M61").Value = -869.3635999999999
$ws.Range("N61").Value = -1474
$ws.Range("H63").Value = 2314.8572
$ws.Range("I63").Value = 2315.6
$ws.Range("J63").Value = 2300
$ws.Range("K63").Value = 2315.6
$ws.Range("L63").Value = 2300
$ws.Range("M63").Value = -1629.6
$ws.Range("N63").Value = -3672
$ws.Range("H66").Value = 2314.8572
$ws.Range("I66").Value = 2315.6
$ws.Range("J66").Value = 2300
$ws.Range("K66").Value = 11578
$ws.Range("L66").Value = 11500
$ws.Range("M66").Value = -8146
$ws.Range("N66").Value = -18364
$ws.Range("H74").Value = 2562.818
$ws.Range("I74").Value = 2687.889
$ws.Range("K74").Value = 2687.889
$ws.Range("M74").Value = -1813.889
$ws.Range("H77").Value = 2562.818
$ws.Range("I77").Value = 2687.889
$ws.Range("K77").Value = 13439.445
$ws.Range("M77").Value = -9071.445
$ws.Range("H122").Value = 1769.2307
$ws.Range("I122").Value = 1769.2307
$ws.Range("K122").Value = 5307.6921
$ws.Range("M122").Value = -2857.6921
$ws.Range("H132").Value = 1308403
$ws.Range("I132").Value = 912.8684
$ws.Range("K132").Value = 2738.6052
$ws.Range("M132").Value = -208.6052
$ws.Range("H136").Value = 1077.973
$ws.Range("I136").Value = 1081.3636
$ws.Range("J136").Value = 1050
$ws.Range("K136").Value = 3244.0908
$ws.Range("L136").Value = 3150
$ws.Range("M136").Value = -694.0907999999999
$ws.Range("N136").Value = -8250

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 329.66666
$ws.Range("I12").Value = 329.66666
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 329.66666
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -161.66666
$ws.Range("N12").ClearContents()
$ws.Range("H17").Value = 2586.3333
$ws.Range("J17").Value = 2586.3333
$ws.Range("L17").Value = 2586.3333
$ws.Range("N17").Value = -2930.3333

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 4011
$ws.Range("I25").Value = 4011
$ws.Range("K25").Value = 4011
$ws.Range("M25").Value = -3837
$ws.Range("H31").Value = 1118.3125
$ws.Range("I31").Value = 798.96295
$ws.Range("K31").Value = 798.96295
$ws.Range("M31").Value = -503.96295
$ws.Range("H34").Value = 1118.3125
$ws.Range("I34").Value = 798.96295
$ws.Range("K34").Value = 798.96295
$ws.Range("M34").Value = -596.96295
$ws.Range("H62").Value = 3018.85
$ws.Range("I62").Value = 2133.1667
$ws.Range("J62").Value = 4347.375
$ws.Range("K62").Value = 2133.1667
$ws.Range("L62").Value = 4347.375
$ws.Range("M62").Value = -1509.1667
$ws.Range("N62").Value = -5595.375
$ws.Range("H65").Value = 3018.85
$ws.Range("I65").Value = 2133.1667
$ws.Range("J65").Value = 4347.375
$ws.Range("K65").Value = 10665.8335
$ws.Range("L65").Value = 21736.875
$ws.Range("M65").Value = -7545.833500000001
$ws.Range("N65").Value = -27976.875

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 23813586
$ws.Range("I5").Value = 31746404
$ws.Range("J5").Value = 15129.286
$ws.Range("K5").Value = 95239212
$ws.Range("L5").Value = 45387.858
$ws.Range("M5").Value = -95239100
$ws.Range("N5").Value = -45611.858
$ws.Range("H122").Value = 7250313
$ws.Range("I122").Value = 29412078
$ws.Range("J122").Value = 5120.423
$ws.Range("K122").Value = 264708702
$ws.Range("L122").Value = 46083.807
$ws.Range("M122").Value = -264706252
$ws.Range("N122").Value = -50983.807
$ws.Range("H131").Value = 884.47
$ws.Range("J131").Value = 919.65216
$ws.Range("L131").Value = 2758.95648
$ws.Range("N131").Value = -12838.95648
$ws.Range("H135").Value = 23813586
$ws.Range("I135").Value = 31746404
$ws.Range("J135").Value = 15129.286
$ws.Range("K135").Value = 285717636
$ws.Range("L135").Value = 136163.574
$ws.Range("M135").Value = -285715101
$ws.Range("N135").Value = -141233.574

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 671.75
$ws.Range("I41").Value = 543.3333
$ws.Range("J41").Value = 1057
$ws.Range("K41").Value = 543.3333
$ws.Range("L41").Value = 1057
$ws.Range("M41").Value = -188.3333
$ws.Range("N41").Value = -1767

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H30").Value = 1620
$ws.Range("I30").Value = 484.4
$ws.Range("J30").Value = 4459
$ws.Range("K30").Value = 484.4
$ws.Range("L30").Value = 4459
$ws.Range("M30").Value = -376.4
$ws.Range("N30").Value = -4675
$ws.Range("H46").Value = 4632025
$ws.Range("I46").Value = 13889415
$ws.Range("J46").Value = 3330
$ws.Range("K46").Value = 13889415
$ws.Range("L46").Value = 3330
$ws.Range("M46").Value = -13889227
$ws.Range("N46").Value = -3706
$ws.Range("H61").Value = 1304.0385
$ws.Range("I61").Value = 1293.4166
$ws.Range("J61").Value = 1313.1428
$ws.Range("K61").Value = 1293.4166
$ws.Range("L61").Value = 1313.1428
$ws.Range("M61").Value = -1091.4166
$ws.Range("N61").Value = -1717.1428
$ws.Range("H113").Value = 1304.0385
$ws.Range("I113").Value = 1293.4166
$ws.Range("J113").Value = 1313.1428
$ws.Range("K113").Value = 1293.4166
$ws.Range("L113").Value = 1313.1428
$ws.Range("M113").Value = 876.5834
$ws.Range("N113").Value = -5653.1428
$ws.Range("H132").Value = 22734290
$ws.Range("I132").Value = 50002520
$ws.Range("K132").Value = 150007560
$ws.Range("M132").Value = -150005030

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H141").Value = 44285.57
$ws.Range("I141").Value = 18500
$ws.Range("J141").Value = 54599.8
$ws.Range("K141").Value = 18500
$ws.Range("L141").Value = 54599.8
$ws.Range("M141").Value = -13320
$ws.Range("N141").Value = -64959.8
